$p = $ppt.ActivePresentation

# --- Reorder slides -------------------------------------------------
# Before: [Triggering/google+pic] [CAEN] [empty Date:/Source:] [empty] [empty]
# After:  [empty] [empty] [Triggering/google+pic] [CAEN] [empty Date:/Source:]
# i.e. the last two (empty) slides move to the front, keeping the
# relative order of the rest intact.

# Move the slide currently last (old slide 5, empty) to position 1.
$p.Slides.Item($p.Slides.Count).MoveTo(1)

# Move the slide now last again (old slide 4, empty) to position 2.
$p.Slides.Item($p.Slides.Count).MoveTo(2)

# --- Add the new question textbox on the slide now at position 1 ----
$s1 = $p.Slides.Item(1)
$tb = $s1.Shapes.AddTextbox(1, 22.588188976377953, 21.176456692913387, 607.7647244094488, 29.081259842519685)
$tb.Name = "TextBox 1"
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0
$tb.Line.Visible = 0
$tb.TextFrame.TextRange.Text = "Q1: What is an event? What is an entry? What is multiplicity?"

# --- Trim the "Source" citation on the slide now at position 4 ------
$s4 = $p.Slides.Item(4)
$sh = $s4.Shapes.Item(1)
$sourcePara = $sh.TextFrame.TextRange.Paragraphs(2, 1)
$sourcePara.Text = "Source:"
